# Apply cryptos list update (Mon Nov 13 11:23:04 UTC 2023 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.926.89"
$ws.Range("E2").Value = "  -0.55%  "

# Row 3
$ws.Range("D3").Value = "2.049.34"
$ws.Range("E3").Value = "  -0.13%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").Value = "'246.43"
$ws.Range("E5").Value = "  -1.20%  "

# Row 6
$ws.Range("E6").Value = "  -2.10%  "

# Row 7
$ws.Range("D7").Value = "'57.87"
$ws.Range("E7").Value = "  -3.73%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").Value = "'0.372"
$ws.Range("E9").Value = "  -4.27%  "

# Row 10
$ws.Range("D10").Value = "'0.0780"
$ws.Range("E10").Value = "  -1.76%  "

# Row 11
$ws.Range("D11").Value = "'0.110"
$ws.Range("E11").Value = "  +2.11%  "

# Row 12
$ws.Range("D12").Value = "'15.26"
$ws.Range("E12").Value = "  -5.29%  "

# Row 13
$ws.Range("D13").Value = "'0.873"
$ws.Range("E13").Value = "  +4.67%  "

# Row 14
$ws.Range("D14").Value = "2.355.24"
$ws.Range("E14").Value = "  +0.21%  "

# Row 15
$ws.Range("D15").Value = "'5.63"
$ws.Range("E15").Value = "  -2.95%  "

# Row 16
$ws.Range("D16").Value = "2.069.77"
$ws.Range("E16").Value = "  +0.80%  "

# Row 17
$ws.Range("D17").Value = "'17.91"
$ws.Range("E17").Value = "  -1.84%  "

# Row 18
$ws.Range("D18").Value = "36.853.46"
$ws.Range("E18").Value = "  -0.79%  "

# Row 19
$ws.Range("D19").Value = "'73.60"
$ws.Range("E19").Value = "  -3.41%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0891"
$ws.Range("E20").Value = "  -1.76%  "

# Row 21
$ws.Range("D21").Value = "'5.40"
$ws.Range("E21").Value = "  -0.50%  "

# Row 22
$ws.Range("D22").Value = "'235.78"
$ws.Range("E22").Value = "  -1.17%  "

# Row 23
$ws.Range("E23").Value = "  +0.02%  "

# Row 24
$ws.Range("E24").Value = "  +1.33%  "

# Row 25
$ws.Range("D25").Value = "'10.32"
$ws.Range("E25").Value = "  +9.15%  "

# Row 26
$ws.Range("D26").Value = "'2.20"
$ws.Range("E26").Value = "  -0.82%  "

# Row 27
$ws.Range("D27").Value = "'168.61"
$ws.Range("E27").Value = "  -0.31%  "

# Row 28
$ws.Range("E28").Value = "  -1.30%  "

# Row 29
$ws.Range("D29").Value = "'5.48"
$ws.Range("E29").Value = "  +13.42%  "

# Row 30
$ws.Range("D30").Value = "'0.123"
$ws.Range("E30").Value = "  -2.11%  "

# Row 31
$ws.Range("D31").Value = "'1.10"
$ws.Range("E31").Value = "  -3.59%  "

# Row 32
$ws.Range("D32").Value = "'4.70"
$ws.Range("E32").Value = "  +1.54%  "

# Row 33
$ws.Range("D33").Value = "'0.0614"
$ws.Range("E33").Value = "  -2.62%  "

# Row 34
$ws.Range("D34").Value = "'2.35"
$ws.Range("E34").Value = "  +4.92%  "

# Row 35
$ws.Range("E35").Value = "  -0.01%  "

# Row 36
$ws.Range("E36").Value = "  +4.16%  "

# Row 37
$ws.Range("D37").Value = "'0.0820"
$ws.Range("E37").Value = "  -8.00%  "

# Row 38
$ws.Range("E38").Value = "  -1.96%  "

# Row 39
$ws.Range("D39").Value = "'5.15"
$ws.Range("E39").Value = "  -0.51%  "

# Row 40
$ws.Range("E40").Value = "  -4.73%  "

# Row 41
$ws.Range("D41").Value = "'0.0223"
$ws.Range("E41").Value = "  -0.61%  "

# Row 42
$ws.Range("D42").Value = "'1.14"
$ws.Range("E42").Value = "  +0.32%  "

# Row 43
$ws.Range("D43").Value = "'0.0949"
$ws.Range("E43").Value = "  -12.68%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'96.90"
$ws.Range("E44").Value = "  -0.49%  "

# Row 45
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'16.91"
$ws.Range("E45").Value = "  -4.38%  "

# Row 46
$ws.Range("D46").Value = "1.302.39"
$ws.Range("E46").Value = "  +0.58%  "

# Row 47
$ws.Range("D47").Value = "'2.36"
$ws.Range("E47").Value = "  -6.15%  "

# Row 48
$ws.Range("D48").Value = "'2.86"
$ws.Range("E48").Value = "  -0.62%  "

# Row 49
$ws.Range("D49").Value = "'6.75"
$ws.Range("E49").Value = "  -1.40%  "

# Row 50
$ws.Range("D50").Value = "2.237.57"
$ws.Range("E50").Value = "  -0.22%  "

# Row 51
$ws.Range("D51").Value = "'45.14"
$ws.Range("E51").Value = "  +1.71%  "
